# adding Program related test data
$wb = $excel.ActiveWorkbook

# --- ProgramPage: add the new test-data table -----------------------------
$ws = $wb.Worksheets.Item("ProgramPage")

# Write header + data row. Order matters: it drives the shared-string table
# indices so they line up with the target workbook (Programname,
# Programdescription, Automation Testing, CoreJavaProgrammingLeven).
$ws.Range("A1").Value = "Programname"
$ws.Range("B1").Value = "Programdescription"
$ws.Range("B2").Value = "Automation Testing"
$ws.Range("A2").Value = "CoreJavaProgrammingLeven"

# Column widths for the new table.
$ws.Columns.Item(1).ColumnWidth = 26.666666666666668
$ws.Columns.Item(2).ColumnWidth = 29.5
$ws.Columns.Item(3).ColumnWidth = 32

# Make ProgramPage the active/selected sheet (mirrors the saved workbook
# view: activeTab + tabSelected move from LoginPage to ProgramPage), then
# leave the cursor on B18 as in the target file.
$ws.Activate() | Out-Null
$ws.Range("B18").Select() | Out-Null

# --- LoginPage: minor column width touch-ups -------------------------------
$login = $wb.Worksheets.Item("LoginPage")
$login.Columns.Item(1).ColumnWidth = 31.666666666666668
$login.Columns.Item(2).ColumnWidth = 25.333333333333332
$login.Columns.Item(4).ColumnWidth = 15.333333333333334
$login.Columns.Item(5).ColumnWidth = 19.333333333333332
